# Update Brechas_Ingresos_Region_7: rename Sexo values (Masculino->Hombre,
# Femenino->Mujer) and refresh the saved cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Data rows 2..65: column G holds the "Sexo" value, alternating
# Masculino/Femenino -> Hombre/Mujer for every record.
for ($r = 2; $r -le 65; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G
    if ($cell.Value2 -eq "Masculino") {
        $cell.Value = "Hombre"
    } elseif ($cell.Value2 -eq "Femenino") {
        $cell.Value = "Mujer"
    }
}

# Restore the sheet's saved selection to a single cell.
[void]$ws.Range("Q58").Select()
